# Append two new evaluation result rows (13 and 14) to the log_evaluations
# sheet, for the llama3:8b-text-q5_K_M and llama3:8b-instruct-q5_K_M runs
# evaluated against llama3:70b.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A .. AB (28 columns), one row per new record.
$newRows = @(
    @(
        "llama3:8b-text-q5_K_M", "llama3:70b", 42, 200, 8743.66, 306.4847, 5,
        "llama3_8b_text_q5_K_M_llama3_70b_42_200_val.txt", 617.78, 1.25,
        "llama3_8b_text_q5_K_M_llama3_70b_42_200_test.txt", 0, 4, 321.31, 2.5,
        "llama3_8b_text_q5_K_M_llama3_70b_42_200_val_fewshot.txt", 640.7, 0,
        "llama3_8b_text_q5_K_M_llama3_70b_42_200_test_fewshot.txt", 5888.22, 2, 2, 321.3, 5,
        "llama3_8b_text_q5_K_M_llama3_70b_42_200_val_bootstrap.txt", 647.86, 0,
        "llama3_8b_text_q5_K_M_llama3_70b_42_200_test_bootstrap.txt"
    ),
    @(
        "llama3:8b-instruct-q5_K_M", "llama3:70b", 42, 200, 8123.89, 273.8765, 10,
        "llama3_8b_instruct_q5_K_M_llama3_70b_42_200_val.txt", 484.3, 15,
        "llama3_8b_instruct_q5_K_M_llama3_70b_42_200_test.txt", 0, 4, 322.1, 0,
        "llama3_8b_instruct_q5_K_M_llama3_70b_42_200_val_fewshot.txt", 645.65, 0,
        "llama3_8b_instruct_q5_K_M_llama3_70b_42_200_test_fewshot.txt", 5625.11, 2, 2, 256.75, 17.5,
        "llama3_8b_instruct_q5_K_M_llama3_70b_42_200_val_bootstrap.txt", 516.1, 8.75,
        "llama3_8b_instruct_q5_K_M_llama3_70b_42_200_test_bootstrap.txt"
    )
)

$startRow = 13
$numCols = 28

# Build a true 2-D array so it can be assigned to a multi-cell Range in one shot,
# the same way Excel COM interop expects for Range.Value bulk writes.
$arr = New-Object 'object[,]' $newRows.Count, $numCols
for ($r = 0; $r -lt $newRows.Count; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $arr[$r, $c] = $newRows[$r][$c]
    }
}

$topLeft = $ws.Cells.Item($startRow, 1)
$bottomRight = $ws.Cells.Item($startRow + $newRows.Count - 1, $numCols)
$targetRange = $ws.Range($topLeft, $bottomRight)
$targetRange.Value = $arr

Write-Host "Added $($newRows.Count) rows starting at row $startRow"
